$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null

$ws.Range("H32").Value = 5927.1
$ws.Range("I32").Value = 6098
$ws.Range("J32").Value = 5813.1665
$ws.Range("K32").Value = 6098
$ws.Range("L32").Value = 5813.1665
$ws.Range("M32").Value = -5772
$ws.Range("N32").Value = -6465.1665

$ws.Range("H51").Value = 10078.842
$ws.Range("I51").Value = 9166.666999999999
$ws.Range("J51").Value = 10249.875
$ws.Range("K51").Value = 9166.666999999999
$ws.Range("L51").Value = 10249.875
$ws.Range("M51").Value = -8682.666999999999
$ws.Range("N51").Value = -11217.875

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

$ws.Range("H111").Value = 506
$ws.Range("I111").Value = 341.33334
$ws.Range("K111").Value = 1024.00002
$ws.Range("M111").Value = 2042.99998

$ws.Range("H129").Value = 1588
$ws.Range("I129").Value = 1334.4
$ws.Range("J129").Value = 2222
$ws.Range("K129").Value = 4003.2
$ws.Range("L129").Value = 6666
$ws.Range("M129").Value = 996.7999999999997
$ws.Range("N129").Value = -16666

$ws.Range("H136").Value = 400000
$ws.Range("J136").Value = 400000
$ws.Range("L136").Value = 400000
$ws.Range("N136").Value = -410200

$ws.Range("H138").Value = 29413770
$ws.Range("J138").Value = 55557580
$ws.Range("L138").Value = 166672740
$ws.Range("N138").Value = -166683020

$ws.Range("H141").Value = 2718.625
$ws.Range("I141").Value = 2392.9048
$ws.Range("K141").Value = 7178.714399999999
$ws.Range("M141").Value = -1998.714399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3921.05
$ws.Range("I74").Value = 3301.1667
$ws.Range("K74").Value = 3301.1667
$ws.Range("M74").Value = -2427.1667

$ws.Range("H77").Value = 3921.05
$ws.Range("I77").Value = 3301.1667
$ws.Range("K77").Value = 16505.8335
$ws.Range("M77").Value = -12137.8335

$ws.Range("H81").Value = 1000000000
$ws.Range("J81").Value = 1000000000
$ws.Range("L81").Value = 1000000000
$ws.Range("N81").Value = -1000001996

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = $null

$ws.Range("H84").Value = 1000000000
$ws.Range("J84").Value = 1000000000
$ws.Range("L84").Value = 3000000000
$ws.Range("N84").Value = -3000009984

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4906.276
$ws.Range("J31").Value = 4536.857
$ws.Range("L31").Value = 4536.857
$ws.Range("N31").Value = -5126.857

$ws.Range("H34").Value = 4906.276
$ws.Range("J34").Value = 4536.857
$ws.Range("L34").Value = 4536.857
$ws.Range("N34").Value = -4940.857

$ws.Range("H35").Value = 3853.6667
$ws.Range("I35").Value = 2669
$ws.Range("J35").Value = 8000
$ws.Range("K35").Value = 2669
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = -2375
$ws.Range("N35").Value = -8588

$ws.Range("H58").Value = 7835.1763
$ws.Range("I58").Value = 6939.8
$ws.Range("J58").Value = 8208.25
$ws.Range("K58").Value = 6939.8
$ws.Range("L58").Value = 8208.25
$ws.Range("M58").Value = -6736.8
$ws.Range("N58").Value = -8614.25

$ws.Range("H99").Value = 6470.8
$ws.Range("I99").Value = 6088.5
$ws.Range("K99").Value = 6088.5
$ws.Range("M99").Value = -4590.5

$ws.Range("H126").Value = 6470.8
$ws.Range("I126").Value = 6088.5
$ws.Range("K126").Value = 18265.5
$ws.Range("M126").Value = -15795.5

$ws.Range("H132").Value = 2223.25
$ws.Range("I132").Value = 2271.842
$ws.Range("K132").Value = 6815.526
$ws.Range("M132").Value = -4285.526

$ws.Range("H133").Value = 42442
$ws.Range("I133").Value = 43000
$ws.Range("K133").Value = 43000
$ws.Range("M133").Value = -40470

$ws.Range("H134").Value = 5240.033
$ws.Range("I134").Value = 4050.5264
$ws.Range("K134").Value = 12151.5792
$ws.Range("M134").Value = -9616.5792

$ws.Range("H136").Value = 7835.1763
$ws.Range("I136").Value = 6939.8
$ws.Range("J136").Value = 8208.25
$ws.Range("K136").Value = 20819.4
$ws.Range("L136").Value = 24624.75
$ws.Range("M136").Value = -18269.4
$ws.Range("N136").Value = -29724.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1793.7
$ws.Range("J12").Value = 2037
$ws.Range("L12").Value = 6111
$ws.Range("N12").Value = -6457

$ws.Range("H24").Value = 1228
$ws.Range("I24").Value = 1228
$ws.Range("K24").Value = 3684
$ws.Range("M24").Value = -3454

$ws.Range("H98").Value = 288.25
$ws.Range("I98").Value = 284.33334
$ws.Range("K98").Value = 853.0000200000001
$ws.Range("M98").Value = 644.9999799999999

$ws.Range("H104").Value = 3761.5
$ws.Range("I104").Value = 3333.8
$ws.Range("J104").Value = 5900
$ws.Range("K104").Value = 10001.4
$ws.Range("L104").Value = 17700
$ws.Range("M104").Value = -7380.400000000001
$ws.Range("N104").Value = -22942

$ws.Range("H122").Value = 350
$ws.Range("I122").Value = 350
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3150
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -700
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5132
$ws.Range("I132").Value = 4591.1113
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 13773.3339
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -11243.3339
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3166.3333
$ws.Range("I7").Value = 2499.5
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 2499.5
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -2387.5
$ws.Range("N7").Value = -4724

$ws.Range("H46").Value = 10418.904
$ws.Range("I46").Value = 3314.75
$ws.Range("K46").Value = 3314.75
$ws.Range("M46").Value = -3126.75

$ws.Range("H93").Value = 2475
$ws.Range("I93").Value = 2450
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 2450
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -1202
$ws.Range("N93").Value = -4996

$ws.Range("H100").Value = 3890.7896
$ws.Range("I100").Value = 2132.75
$ws.Range("K100").Value = 2132.75
$ws.Range("M100").Value = -1591.75

$ws.Range("H126").Value = 3166.3333
$ws.Range("I126").Value = 2499.5
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 7498.5
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -5028.5
$ws.Range("N126").Value = -18440

$ws.Range("H132").Value = 7460.0425
$ws.Range("I132").Value = 7160.775
$ws.Range("K132").Value = 21482.325
$ws.Range("M132").Value = -18952.325

$ws.Range("H136").Value = 4382.913
$ws.Range("I136").Value = 3637.4443
$ws.Range("K136").Value = 10912.3329
$ws.Range("M136").Value = -8362.332900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 20000
$ws.Range("J30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("N30").Value = -20214

$ws.Range("H86").Value = 88000
$ws.Range("J86").Value = 88000
$ws.Range("L86").Value = 88000
$ws.Range("N86").Value = -90246

$ws.Range("H89").Value = 88000
$ws.Range("J89").Value = 88000
$ws.Range("L89").Value = 440000
$ws.Range("N89").Value = -451232

$ws.Range("H122").Value = 2516.4736
$ws.Range("I122").Value = 1518.4117
$ws.Range("K122").Value = 4555.2351
$ws.Range("M122").Value = -2105.2351

$ws.Range("H126").Value = 1576.3684
$ws.Range("I126").Value = 1576.3684
$ws.Range("K126").Value = 4729.1052
$ws.Range("M126").Value = -2259.1052

$ws.Range("H132").Value = 2533.3704
$ws.Range("I132").Value = 2164.06
$ws.Range("K132").Value = 6492.18
$ws.Range("M132").Value = -3962.18

$ws.Range("H136").Value = 6930.381
$ws.Range("I136").Value = 6036.643
$ws.Range("J136").Value = 8717.857
$ws.Range("K136").Value = 18109.929
$ws.Range("L136").Value = 26153.571
$ws.Range("M136").Value = -15559.929
$ws.Range("N136").Value = -31253.571
